# Edit: insert two new weekly price rows for "Lane Late" orange variety
# at the top of the Vega Monumental Concepción block (rows 336-337),
# pushing the existing rows 336-351 down to 338-353.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 336 (shifts old 336.. down by 2)
$ws.Range("A336:T337").EntireRow.Insert()

# --- New row 336 ---
$ws.Cells.Item(336, 1).Value = 11
$ws.Cells.Item(336, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(336, 3).Value = "Bíobío"
$ws.Cells.Item(336, 4).Value = 44826
$ws.Cells.Item(336, 5).Value = 8
$ws.Cells.Item(336, 6).Value = "Fruta"
$ws.Cells.Item(336, 7).Value = 100102
$ws.Cells.Item(336, 8).Value = "Cítricos"
$ws.Cells.Item(336, 9).Value = 100102005
$ws.Cells.Item(336, 10).Value = "Naranja"
$ws.Cells.Item(336, 11).Value = "Lane Late"
$ws.Cells.Item(336, 12).Value = "Primera"
$ws.Cells.Item(336, 13).Value = 150
$ws.Cells.Item(336, 14).Value = 7000
$ws.Cells.Item(336, 15).Value = 7000
$ws.Cells.Item(336, 16).Value = 7000
$ws.Cells.Item(336, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(336, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(336, 19).Value = 467
$ws.Cells.Item(336, 20).Value = 15

# --- New row 337 ---
$ws.Cells.Item(337, 1).Value = 11
$ws.Cells.Item(337, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(337, 3).Value = "Bíobío"
$ws.Cells.Item(337, 4).Value = 44826
$ws.Cells.Item(337, 5).Value = 8
$ws.Cells.Item(337, 6).Value = "Fruta"
$ws.Cells.Item(337, 7).Value = 100102
$ws.Cells.Item(337, 8).Value = "Cítricos"
$ws.Cells.Item(337, 9).Value = 100102005
$ws.Cells.Item(337, 10).Value = "Naranja"
$ws.Cells.Item(337, 11).Value = "Lane Late"
$ws.Cells.Item(337, 12).Value = "Segunda"
$ws.Cells.Item(337, 13).Value = 200
$ws.Cells.Item(337, 14).Value = 6000
$ws.Cells.Item(337, 15).Value = 6000
$ws.Cells.Item(337, 16).Value = 6000
$ws.Cells.Item(337, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(337, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(337, 19).Value = 400
$ws.Cells.Item(337, 20).Value = 15

# Ensure the date cells use the same date/time number format as the rest of column D
$ws.Range("D336:D337").NumberFormat = $ws.Cells.Item(338, 4).NumberFormat
